$d = $word.ActiveDocument
$p2 = $d.Paragraphs.Item(2)
$pr = $p2.Range
$searchText = "Invalid block: Unexpected tag EOF missing [ENDFOR] while parsing m:for v | self.eClassifiers"

$searchRange = $d.Range($pr.Start, $pr.End)
$count = 0
$found = $searchRange.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $count = $count + 1
    $searchRange.Font.Bold = 1
    $endPos = $searchRange.End
    $searchRange.Start = $endPos
    $searchRange.End = $pr.End
    $found = $searchRange.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}
Write-Host ("Total matches updated: " + $count)
